$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Some of the new "article code" values are purely numeric strings
# (e.g. "310405154321"). Assigning such a string directly to .Value would
# make Excel auto-detect it as a number. Force a Text number format first so
# the value is stored as text, matching the original column's text data.
# (The row-wide format paste further below restores the correct style/fill
#  while leaving the already-stored text value untouched.)
# ---------------------------------------------------------------------------
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "310405154321"

$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "350211164221"

$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "350211064221"

$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "10350202464221"

$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "350103454326"

# ---------------------------------------------------------------------------
# Apply the highlighted "match found" style (same as rows 5 & 6, style s="2",
# light-green fill) to rows 2, 3, 4, 7, 8, 9, 10, 11, 12 by copying formats
# from an already-styled row. This reuses the existing style instead of
# creating new style entries.
# ---------------------------------------------------------------------------
$ws.Range("A5:I5").Copy()
$ws.Range("A2:I2").PasteSpecial(-4122)
$ws.Range("A3:I3").PasteSpecial(-4122)
$ws.Range("A4:I4").PasteSpecial(-4122)
$ws.Range("A7:I7").PasteSpecial(-4122)
$ws.Range("A8:I8").PasteSpecial(-4122)
$ws.Range("A9:I9").PasteSpecial(-4122)
$ws.Range("A10:I10").PasteSpecial(-4122)
$ws.Range("A11:I11").PasteSpecial(-4122)
$ws.Range("A12:I12").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Row 2 - Барабан тормозной МАЗ: updated match percentage
# ---------------------------------------------------------------------------
$ws.Range("G2").Value = 100

# ---------------------------------------------------------------------------
# Row 3 - Болт короткий МАЗ: now matched to a different nomenclature item
# ---------------------------------------------------------------------------
$ws.Range("B3").Value = "54321-3104051"
$ws.Range("D3").Value = "54321-3104051-01-СПЕЦМАШ"
$ws.Range("E3").Value = "Болт ступицы (Еврошпилька М22х1,5 L-82 мм) МАЗ, ТЕФЛОН (упак. 10 шт)"
$ws.Range("F3").Value = "01310405154321СПЕЦМАШ"
$ws.Range("G3").Value = 100
$ws.Range("H3").Value = 141

# ---------------------------------------------------------------------------
# Row 4 - Втулка цапфы МАЗ: updated match percentage
# ---------------------------------------------------------------------------
$ws.Range("G4").Value = 100

# ---------------------------------------------------------------------------
# Row 7 - Кулак разжимной (лев): now matched to a different nomenclature item
# ---------------------------------------------------------------------------
$ws.Range("B7").Value = "64221-3502111"
$ws.Range("D7").Value = "64221-3502111-10-СПЕЦМАШ"
$ws.Range("E7").Value = "Кулак разжимной (L=546 мм, левый, под рычаг с широким шлицем) "
$ws.Range("F7").Value = "10350211164221СПЕЦМАШ"
$ws.Range("G7").Value = 100
$ws.Range("H7").Value = 2775

# ---------------------------------------------------------------------------
# Row 8 - Кулак разжимной (прав): now matched to a different nomenclature item
# ---------------------------------------------------------------------------
$ws.Range("B8").Value = "64221-3502110"
$ws.Range("D8").Value = "64221-3502110-10-СПЕЦМАШ"
$ws.Range("E8").Value = "Кулак разжимной (L=546 мм, правый, под рычаг с широким шлицем) "
$ws.Range("F8").Value = "10350211064221СПЕЦМАШ"
$ws.Range("G8").Value = 100
$ws.Range("H8").Value = 2775

# ---------------------------------------------------------------------------
# Row 9 - Опора кулака разжимного: now matched to a different nomenclature item
# ---------------------------------------------------------------------------
$ws.Range("B9").Value = "64221-3502024-10"
$ws.Range("D9").Value = "64221-3502024-СПЕЦМАШ"
$ws.Range("E9").Value = "Опора разжимного кулака в сборе (с ШС-40) МАЗ-5336, 6422, 103, 104, 105"
$ws.Range("F9").Value = "350202464221СПЕЦМАШ"
$ws.Range("G9").Value = 100
$ws.Range("H9").Value = 1595

# ---------------------------------------------------------------------------
# Row 10 - Пружина 54326: now matched to a different nomenclature item
# ---------------------------------------------------------------------------
$ws.Range("B10").Value = "54326-3501034"
$ws.Range("D10").Value = "54326-3501034-СПЕЦМАШ"
$ws.Range("E10").Value = "Пружина полуприцепа колодок стяжная (двойная) (упак. 10 шт) для МАЗ"
$ws.Range("F10").Value = "350103454326СПЕЦМАШ"
$ws.Range("G10").Value = 100
$ws.Range("H10").Value = 123

# ---------------------------------------------------------------------------
# Row 11 - Тяга рулевая поперечная: updated match percentage
# ---------------------------------------------------------------------------
$ws.Range("G11").Value = 100

# ---------------------------------------------------------------------------
# Row 12 - Тяга рулевая продольная: now matched to a different nomenclature item
# ---------------------------------------------------------------------------
$ws.Range("D12").Value = "5551-3003010-СПЕЦМАШ"
$ws.Range("E12").Value = "Тяга рулевая продольная МАЗ-533702, 543202, 54323, 5551, 555102, 5516, 03, 05, 6422, 64229, L=754 мм"
$ws.Range("F12").Value = "30030105551СПЕЦМАШ"
$ws.Range("G12").Value = 100
$ws.Range("H12").Value = 9795
